# Add data to excel import:
#  - new "data" worksheet (after "forms") holding the "address-data" record
#  - "forms" sheet header relabelled "Form Id (ident)" / column A width added
#  - selections / active sheet restored to match

$wb = $excel.ActiveWorkbook
$formsSheet = $wb.Worksheets.Item("forms")

# ---------------------------------------------------------------------
# 1. forms sheet: header text + new column A width
# ---------------------------------------------------------------------
$formsSheet.Range("A1").Value = "Form Id (ident)"
$formsSheet.Columns.Item(1).ColumnWidth = 16.0

# ---------------------------------------------------------------------
# 2. add the new "data" worksheet right after "forms"
# ---------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $formsSheet)
$dataSheet.Name = "data"

# copy the header/content formatting (fills + fonts) from the forms sheet
# so the new sheet reuses the same style records instead of minting new ones
$formsSheet.Range("A1:B1").Copy() | Out-Null
$dataSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$formsSheet.Range("B2").Copy() | Out-Null
$dataSheet.Range("B2").PasteSpecial(-4122) | Out-Null

# column widths matching the forms sheet
$dataSheet.Columns.Item(1).ColumnWidth = 16.0
$dataSheet.Columns.Item(2).ColumnWidth = 91.85

# values
$dataSheet.Range("A1").Value = "Data Id (ident)"
$dataSheet.Range("B1").Value = "Content"
$dataSheet.Range("A2").Value = "address-data"
$dataSheet.Range("B2").Value = '{"ident":"address-data","structure":{"value":{"street":{"DataString":{"value":"Sonnenweg"}},"number":{"DataString":{"value":"23a"}},"postcode":{"DataNumber":{"value":6414}},"city":{"DataString":{"value":"Oberarth"}}}}}'

# match zoom level with the forms sheet, and the saved selection
$dataSheet.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 170
$dataSheet.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. restore forms sheet as the active tab / selection
# ---------------------------------------------------------------------
$formsSheet.Activate() | Out-Null
$formsSheet.Range("B2").Select() | Out-Null

Write-Output "done"
